$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "sdsad"
$ws.Range("D7").Value = "sdsad"
$ws.Range("G8").Value = "sdsa"
$ws.Range("H13").Value = "sdsad"

$ws.Range("G10").Select()
